$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '94.604.74'
$ws.Range("E2").Value = '  -3.59%  '

$ws.Range("D3").Value = '3.427.34'
$ws.Range("E3").Value = '  +1.31%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").Value = '''238.57'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.35%  '

$ws.Range("D6").Value = '''641.78'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.51%  '

$ws.Range("D7").Value = '''1.43'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.49%  '

$ws.Range("D8").Value = '''0.404'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.48%  '

$ws.Range("E9").Value = '  +0.14%  '

$ws.Range("D10").Value = '''0.974'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.03%  '

$ws.Range("D11").Value = '3.427.74'
$ws.Range("E11").Value = '  +1.42%  '

$ws.Range("E12").Value = '  -4.12%  '

$ws.Range("D13").Value = '''41.65'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.61%  '

$ws.Range("D14").Value = '''6.27'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.80%  '

$ws.Range("D15").Value = '94.427.17'
$ws.Range("E15").Value = '  -3.42%  '

$ws.Range("D16").Value = '4.060.10'
$ws.Range("E16").Value = '  +1.28%  '

$ws.Range("D17").Value = '''0.0000252'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.75%  '

$ws.Range("D18").Value = '''8.35'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -9.23%  '

$ws.Range("D19").Value = '3.430.23'
$ws.Range("E19").Value = '  +2.11%  '

$ws.Range("D20").Value = '''17.51'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.74%  '

$ws.Range("D21").Value = '''11.54'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.09%  '

$ws.Range("D22").Value = '''0.488'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.89%  '

$ws.Range("D23").Value = '''500.41'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.80%  '

$ws.Range("D24").Value = '''3.24'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.05%  '

$ws.Range("E25").Value = '  -3.97%  '

$ws.Range("D26").Value = '''6.52'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.23%  '

$ws.Range("D27").Value = '''91.56'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.20%  '

$ws.Range("D28").Value = '3.612.66'
$ws.Range("E28").Value = '  +1.27%  '

$ws.Range("D29").Value = '''11.94'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.98%  '

$ws.Range("D30").Value = '''11.61'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.08%  '

$ws.Range("E31").Value = '  +0.16%  '

$ws.Range("D32").Value = '''2.74'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.54%  '

$ws.Range("D33").Value = '''0.137'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.24%  '

$ws.Range("D34").Value = '''0.179'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.90%  '

$ws.Range("D35").Value = '''0.997'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.39%  '

$ws.Range("D36").Value = '''29.63'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.94%  '

$ws.Range("D37").Value = '''0.552'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.17%  '

$ws.Range("D38").Value = '''544.83'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.52%  '

$ws.Range("D39").Value = '''7.69'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.62%  '

$ws.Range("E40").Value = '  -1.79%  '

$ws.Range("B41").Value = 'USDe'
$ws.Range("C41").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D41").Value = '''1.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.02%  '

$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Value = '''0.151'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.82%  '

$ws.Range("D43").Value = '''0.908'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +7.25%  '

$ws.Range("D44").Value = '''24.07'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.40%  '

$ws.Range("E45").Value = '  +0.61%  '

$ws.Range("E46").Value = '  +0.15%  '

$ws.Range("D47").Value = '''5.65'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.99%  '

$ws.Range("D48").Value = '''2.21'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.94%  '

$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value = '''0.0410'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.49%  '

$ws.Range("B50").Value = 'dogwifhat'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D50").Value = '''3.32'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.28%  '

$ws.Range("D51").Value = '''54.88'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.07%  '
